$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rubric")
$ws.Activate()

# --- Delete column C (values mirrored column B but were otherwise unused) ---
# This shifts D->C, E->D, F->E exactly like the target diff shows (D2/D13 appear
# with the old E2/E13 style, E2/E13 keep their style, F2/F13 disappear).
$ws.Columns("C:C").Delete()

# Column B width becomes a custom 8.5 (character width units as stored in xml).
# ColumnWidth (COM, character units) differs from the stored xml "width" by
# the standard 5/6 padding constant, so 8.5 - 5/6 = 7.666... gives xml width 8.5.
$ws.Columns("B:B").ColumnWidth = 7.666666666666667

# --- Header label: "Possible" -> "Points" ---
$ws.Range("B2").Value = "Points"

# --- Point values revised throughout the rubric ---
$ws.Range("B4").Value = 25
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 12
$ws.Range("B8").Value = 25
$ws.Range("B9").Value = 15
$ws.Range("B10").Value = 12
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 25

$ws.Range("B16").Value = 10
$ws.Range("B17").Value = 20
$ws.Range("B18").Value = 12
$ws.Range("B19").Value = 5
$ws.Range("B20").Value = 10
$ws.Range("B21").Value = 10

$ws.Range("B25").Value = 15
$ws.Range("B26").Value = 30

$ws.Range("B29").Value = 3

# --- Add a thin bottom border under the last item of each section (the row
# immediately above each "Subtotal" line), matching the new border/style. ---
foreach ($addr in @("B12","B21","B26")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = 2
}

# --- Selection / scroll position left at C29 (no frozen top-left override) ---
$ws.Range("C29").Select()
